$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.738708734512329
$ws.Range("B1").Value = 2.795881032943726
$ws.Range("C1").Value = 3.480600118637085
$ws.Range("D1").Value = 1.339206099510193
$ws.Range("E1").Value = 0.8928744196891785
